$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q3" sheet before the existing "2022-Q2" sheet.
#    Copying the existing "2022-Q2" sheet keeps identical formatting
#    (styles, column layout, dimension) for the new tab.
# ------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Copy($oldQ2)
$newQ3 = $wb.Worksheets.Item(2)
$newQ3.Name = "2022-Q3"

# Update the fund rows on the new "2022-Q3" sheet with the Q3 figures
# (fund code / fund name columns stay the same, only the metrics change).
# These columns hold text-formatted numbers (e.g. "45.00"), so force a
# text number format before assigning, otherwise Excel coerces the
# string into a numeric value and the trailing zeroes are lost.
$textCells = $newQ3.Range("D2:G3")
$textCells.NumberFormat = "@"

$newQ3.Range("D2").Value = "2.70"
$newQ3.Range("E2").Value = "45.00"
$newQ3.Range("F2").Value = "2.44"
$newQ3.Range("G2").Value = "0.0659"
$newQ3.Range("H2").Value = 2

$newQ3.Range("D3").Value = "1.73"
$newQ3.Range("E3").Value = "45.00"
$newQ3.Range("F3").Value = "2.44"
$newQ3.Range("G3").Value = "0.0422"
$newQ3.Range("H3").Value = 2

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3
#    above the 2022-Q2 row, pushing everything else down.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Match the formatting of the other data rows (Insert() alone doesn't
# reliably carry over the bordered/bold style used by column A, nor the
# plain style used by columns B-D), by copying it from the row below.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

# Renumber the index column (A) for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Restore the originally-active tab ("2021-Q2", the last sheet) since
# copying a new sheet in shifts the active selection to it.
$wb.Worksheets.Item("2021-Q2").Activate()
